$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column H (Examen) marks that changed from 0 to their real value ---
$ws.Range("H2").Value = 85
$ws.Range("H3").Value = 67
$ws.Range("H4").Value = 50
$ws.Range("H5").Value = 61
$ws.Range("H6").Value = 58
$ws.Range("H7").Value = 58
$ws.Range("H8").Value = 43
$ws.Range("H10").Value = 61
$ws.Range("H11").Value = 30
$ws.Range("H12").Value = 34
$ws.Range("H13").Value = 85
$ws.Range("H14").Value = 82
$ws.Range("H15").Value = 82
$ws.Range("H16").Value = 64
$ws.Range("H17").Value = 34
$ws.Range("H18").Value = 50
$ws.Range("H19").Value = 36
$ws.Range("H21").Value = 85
$ws.Range("H23").Value = 36
$ws.Range("H24").Value = 67
$ws.Range("H26").Value = 30
$ws.Range("H27").Value = 54
$ws.Range("H28").Value = 54
$ws.Range("H29").Value = 54
$ws.Range("H31").Value = 56
$ws.Range("H33").Value = 54
$ws.Range("H34").Value = 21
$ws.Range("H35").Value = 21
$ws.Range("H38").Value = 56
$ws.Range("H40").Value = 67
$ws.Range("H42").Value = 11
$ws.Range("H45").Value = 54
$ws.Range("H46").Value = 54
$ws.Range("H47").Value = 67
$ws.Range("H49").Value = 54
$ws.Range("H50").Value = 44
$ws.Range("H52").Value = 84
$ws.Range("H53").Value = 70
$ws.Range("H55").Value = 51
$ws.Range("H56").Value = 63
$ws.Range("H57").Value = 70
$ws.Range("H58").Value = 31
$ws.Range("H59").Value = 54
$ws.Range("H60").Value = 84
$ws.Range("H61").Value = 51
$ws.Range("H63").Value = 36
$ws.Range("H64").Value = 86
$ws.Range("H65").Value = 55
$ws.Range("H66").Value = 36
$ws.Range("H67").Value = 44
$ws.Range("H68").Value = 86
$ws.Range("H71").Value = 55
$ws.Range("H72").Value = 63
$ws.Range("H73").Value = 80
$ws.Range("H74").Value = 56
$ws.Range("H75").Value = 13
$ws.Range("H76").Value = 95
$ws.Range("H77").Value = 74
$ws.Range("H78").Value = 79
$ws.Range("H79").Value = 79
$ws.Range("H80").Value = 71
$ws.Range("H81").Value = 80
$ws.Range("H82").Value = 56
$ws.Range("H84").Value = 74
$ws.Range("H85").Value = 75
$ws.Range("H86").Value = 71
$ws.Range("H87").Value = 56
$ws.Range("H88").Value = 72
$ws.Range("H89").Value = 85
$ws.Range("H90").Value = 75
$ws.Range("H91").Value = 72
$ws.Range("H92").Value = 76
$ws.Range("H93").Value = 95
$ws.Range("H94").Value = 60
$ws.Range("H95").Value = 80
$ws.Range("H96").Value = 79
$ws.Range("H97").Value = 59
$ws.Range("H98").Value = 62
$ws.Range("H99").Value = 56
$ws.Range("H100").Value = 60
$ws.Range("H101").Value = 76
$ws.Range("H102").Value = 80
$ws.Range("H103").Value = 59
$ws.Range("H104").Value = 79
$ws.Range("H105").Value = 62

# --- Rewrite the "Total" formula (column J) to clamp/round the grade, and
#     add a new helper column L holding the raw (unrounded) grade ---
$ws.Range("J2:J106").FormulaR1C1 = "=MAX(4, ROUND(MIN(10, 1+0.5*RC[-1]/5+0.2*RC[-2]/10+0.06*(RC[-7]/9+RC[-6]/7+RC[-5]/7+RC[-4]/7+RC[-3]/6)),0))"
$ws.Range("L2:L106").FormulaR1C1 = "=MIN(10, 1+0.5*RC[-3]/5+0.2*RC[-4]/10+0.06*(RC[-9]/9+RC[-8]/7+RC[-7]/7+RC[-6]/7+RC[-5]/6))"

# --- Highlight a couple of individual L cells with a yellow fill (as in the source file) ---
$ws.Range("L10").Interior.Color = 65535
$ws.Range("L17").Interior.Color = 65535

# --- Conditional formatting: highlight borderline grades (4.05 - 4.49) in column L ---
$condRange = $ws.Range("L1:L1048576")
$fc = $condRange.FormatConditions.Add(1, 1, "4.05", "4.49")
$fc.Interior.Color = 65535

# --- Worksheet view tweaks ---
$ws.Range("E24").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1

# --- Page setup ---
$ws.PageSetup.Orientation = 1

Write-Host "Edit applied"
